$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.228.65"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").Value = "3.671.01"
$ws.Range("E3").Value = "  -0.44%  "

# Row 4
$ws.Range("E4").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.03"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.62"
$ws.Range("E6").Value = "  +6.81%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.701"
$ws.Range("E9").Value = "  -2.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  -4.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.31"
$ws.Range("E11").Value = "  +1.94%  "

# Row 12
$ws.Range("E12").Value = "  -5.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.25"
$ws.Range("E13").Value = "  -1.17%  "

# Row 14
$ws.Range("D14").Value = "4.253.13"
$ws.Range("E14").Value = "  +0.05%  "

# Row 15
$ws.Range("D15").Value = "3.675.06"
$ws.Range("E15").Value = "  -0.02%  "

# Row 16
$ws.Range("E16").Value = "  +0.44%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.92"
$ws.Range("E17").Value = "  -2.26%  "

# Row 18
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.11"
$ws.Range("E18").Value = "  -0.86%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.63"
$ws.Range("E19").Value = "  -1.10%  "

# Row 20
$ws.Range("D20").Value = "68.109.06"
$ws.Range("E20").Value = "  +0.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "406.34"
$ws.Range("E21").Value = "  -0.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.44"
$ws.Range("E22").Value = "  -1.66%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.55"
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("E24").Value = "  -1.33%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.93"
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("E26").Value = "  -0.84%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.08"
$ws.Range("E27").Value = "  +0.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.72"
$ws.Range("E28").Value = "  -3.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("E29").Value = "  -0.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.03"
$ws.Range("E30").Value = "  -1.78%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.18"
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.34"
$ws.Range("E32").Value = "  -0.96%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "66.93"
$ws.Range("E33").Value = "  +4.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "44.20"
$ws.Range("E34").Value = "  +2.02%  "

# Row 35
$ws.Range("E35").Value = "  +0.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "611.74"
$ws.Range("E36").Value = "  +3.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.09%  "

# Row 38
$ws.Range("E38").Value = "  -1.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.16%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("E40").Value = "  -10.80%  "

# Row 41
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").Value = "  -2.34%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0428"
$ws.Range("E43").Value = "  -1.00%  "

# Row 44
$ws.Range("E44").Value = "  -7.34%  "

# Row 45
$ws.Range("E45").Value = "  +2.18%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.23"
$ws.Range("E46").Value = "  +2.57%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.788.72"
$ws.Range("E47").Value = "  +0.86%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.96"
$ws.Range("E48").Value = "  -2.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.52"
$ws.Range("E49").Value = "  +2.38%  "

# Row 50
$ws.Range("E50").Value = "  -3.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.53"
$ws.Range("E51").Value = "  -11.83%  "
